$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 337.75
$ws.Range("I18").Value = 337.75
$ws.Range("K18").Value = 337.75
$ws.Range("M18").Value = -53.75
$ws.Range("H69").Value = 4000
$ws.Range("I69").Value = 4000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -11126
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 4000
$ws.Range("I72").Value = 4000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 36000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -31632
$ws.Range("N72").ClearContents()
$ws.Range("H76").Value = 5449.3335
$ws.Range("I76").Value = 6174.125
$ws.Range("J76").Value = 3999.75
$ws.Range("K76").Value = 6174.125
$ws.Range("L76").Value = 3999.75
$ws.Range("M76").Value = -5859.125
$ws.Range("N76").Value = -4629.75
$ws.Range("H79").Value = 5449.3335
$ws.Range("I79").Value = 6174.125
$ws.Range("J79").Value = 3999.75
$ws.Range("K79").Value = 6174.125
$ws.Range("L79").Value = 3999.75
$ws.Range("M79").Value = -5082.125
$ws.Range("N79").Value = -6183.75
$ws.Range("H80").Value = 7686.5
$ws.Range("I80").Value = 747.75
$ws.Range("J80").Value = 14625.25
$ws.Range("K80").Value = 2243.25
$ws.Range("L80").Value = 43875.75
$ws.Range("M80").Value = -1245.25
$ws.Range("N80").Value = -45871.75
$ws.Range("H82").Value = 8789.267
$ws.Range("I82").Value = 1530.8182
$ws.Range("J82").Value = 28750
$ws.Range("K82").Value = 4592.4546
$ws.Range("L82").Value = 86250
$ws.Range("M82").Value = -4186.4546
$ws.Range("N82").Value = -87062
$ws.Range("H83").Value = 7686.5
$ws.Range("I83").Value = 747.75
$ws.Range("J83").Value = 14625.25
$ws.Range("K83").Value = 6729.75
$ws.Range("L83").Value = 131627.25
$ws.Range("M83").Value = -1737.75
$ws.Range("N83").Value = -141611.25
$ws.Range("H85").Value = 8789.267
$ws.Range("I85").Value = 1530.8182
$ws.Range("J85").Value = 28750
$ws.Range("K85").Value = 4592.4546
$ws.Range("L85").Value = 86250
$ws.Range("M85").Value = -3188.4546
$ws.Range("N85").Value = -89058
$ws.Range("H88").Value = 2517.25
$ws.Range("I88").Value = 2601
$ws.Range("J88").Value = 2489.3333
$ws.Range("K88").Value = 2601
$ws.Range("L88").Value = 2489.3333
$ws.Range("M88").Value = -2195
$ws.Range("N88").Value = -3301.3333
$ws.Range("H91").Value = 2517.25
$ws.Range("I91").Value = 2601
$ws.Range("J91").Value = 2489.3333
$ws.Range("K91").Value = 2601
$ws.Range("L91").Value = 2489.3333
$ws.Range("M91").Value = -1197
$ws.Range("N91").Value = -5297.3333
$ws.Range("H103").Value = 2118
$ws.Range("I103").Value = 399.33334
$ws.Range("J103").Value = 2762.5
$ws.Range("K103").Value = 1198.00002
$ws.Range("L103").Value = 8287.5
$ws.Range("M103").Value = -612.0000199999999
$ws.Range("N103").Value = -9459.5
$ws.Range("H108").Value = 23833
$ws.Range("J108").Value = 23833
$ws.Range("L108").Value = 23833
$ws.Range("N108").Value = -31513
$ws.Range("H110").Value = 25000
$ws.Range("J110").Value = 25000
$ws.Range("L110").Value = 25000
$ws.Range("N110").Value = -33180
$ws.Range("H129").Value = 1194.6904
$ws.Range("I129").Value = 423.55554
$ws.Range("J129").Value = 1405
$ws.Range("K129").Value = 1270.66662
$ws.Range("L129").Value = 4215
$ws.Range("M129").Value = 3729.33338
$ws.Range("N129").Value = -14215
$ws.Range("H138").Value = 2162.6714
$ws.Range("I138").Value = 1198.1395
$ws.Range("J138").Value = 3698.7778
$ws.Range("K138").Value = 3594.4185
$ws.Range("L138").Value = 11096.3334
$ws.Range("M138").Value = 1545.5815
$ws.Range("N138").Value = -21376.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15913.16
$ws.Range("I32").Value = 16772.465
$ws.Range("K32").Value = 16772.465
$ws.Range("M32").Value = -16485.465
$ws.Range("H34").Value = 21111.111
$ws.Range("J34").Value = 21111.111
$ws.Range("L34").Value = 21111.111
$ws.Range("N34").Value = -21653.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 48799.6
$ws.Range("J20").Value = 48799.6
$ws.Range("L20").Value = 48799.6
$ws.Range("N20").Value = -49271.6
$ws.Range("H30").Value = 48799.6
$ws.Range("J30").Value = 48799.6
$ws.Range("L30").Value = 48799.6
$ws.Range("N30").Value = -48981.6
$ws.Range("H128").Value = 48799.6
$ws.Range("J128").Value = 48799.6
$ws.Range("L128").Value = 48799.6
$ws.Range("N128").Value = -58759.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2363.8708
$ws.Range("I137").Value = 1078.8235
$ws.Range("J137").Value = 3924.2856
$ws.Range("K137").Value = 3236.4705
$ws.Range("L137").Value = 11772.8568
$ws.Range("M137").Value = 1863.5295
$ws.Range("N137").Value = -21972.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 20080660
$ws.Range("I7").Value = 50000150
$ws.Range("J7").Value = 134333.33
$ws.Range("K7").Value = 50000150
$ws.Range("L7").Value = 134333.33
$ws.Range("M7").Value = -50000038
$ws.Range("N7").Value = -134557.33
$ws.Range("H8").Value = 20080660
$ws.Range("I8").Value = 50000150
$ws.Range("J8").Value = 134333.33
$ws.Range("K8").Value = 50000150
$ws.Range("L8").Value = 134333.33
$ws.Range("M8").Value = -50000011
$ws.Range("N8").Value = -134611.33
$ws.Range("H102").Value = 2258.138
$ws.Range("I102").Value = 1812.6842
$ws.Range("J102").Value = 3104.5
$ws.Range("K102").Value = 1812.6842
$ws.Range("L102").Value = 3104.5
$ws.Range("M102").Value = -190.6841999999999
$ws.Range("N102").Value = -6348.5
$ws.Range("H120").Value = 40317
$ws.Range("J120").Value = 40317
$ws.Range("L120").Value = 40317
$ws.Range("N120").Value = -49993

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 17533.166
$ws.Range("J56").Value = 25053.5
$ws.Range("L56").Value = 25053.5
$ws.Range("N56").Value = -26481.5
